$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.243.85'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '1.650.16'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = '''217.34'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").Value = '''0.0629'
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").Value = '1.882.16'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '1.647.68'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  +2.56%  '
$ws.Range("D16").Value = '''67.68'
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("D17").Value = '27.233.56'
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").Value = '''219.99'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").Value = '''6.85'
$ws.Range("E21").Value = '  +2.98%  '
$ws.Range("E22").Value = '  +5.04%  '
$ws.Range("E23").Value = '  +0.68%  '
$ws.Range("D24").Value = '''9.22'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Value = '''147.11'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''7.53'
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").Value = '''0.118'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").Value = '''15.81'
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("D30").Value = '''0.0509'
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").Value = '''3.39'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").Value = '1.263.18'
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  +1.06%  '
$ws.Range("E38").Value = '  +2.76%  '
$ws.Range("D39").Value = '''0.845'
$ws.Range("E39").Value = '  +1.66%  '
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").Value = '''5.44'
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("E43").Value = '  +5.79%  '
$ws.Range("D44").Value = '1.792.26'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '''62.04'
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("D46").Value = '''91.74'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = '''1.60'
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '''7.67'
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").Value = '''0.0971'
$ws.Range("E51").Value = '  -0.46%  '
